$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" "39.044.69"
Set-TextValue $ws "E2" "  -4.87%  "
Set-TextValue $ws "D3" "2.245.84"
Set-TextValue $ws "E3" "  -7.06%  "
Set-TextValue $ws "E4" "  +0.16%  "
Set-TextValue $ws "D5" "292.98"
Set-TextValue $ws "E5" "  -7.24%  "
Set-TextValue $ws "D6" "79.92"
Set-TextValue $ws "E6" "  -9.94%  "
Set-TextValue $ws "D7" "0.504"
Set-TextValue $ws "E7" "  -6.25%  "
Set-TextValue $ws "E8" "  +0.12%  "
Set-TextValue $ws "D9" "0.455"
Set-TextValue $ws "E9" "  -8.33%  "
Set-TextValue $ws "D10" "0.0766"
Set-TextValue $ws "E10" "  -8.16%  "
Set-TextValue $ws "D11" "27.65"
Set-TextValue $ws "E11" "  -12.30%  "
Set-TextValue $ws "D12" "45.47"
Set-TextValue $ws "E12" "  -14.82%  "
Set-TextValue $ws "E13" "  -1.35%  "
Set-TextValue $ws "D14" "2.614.15"
Set-TextValue $ws "E14" "  -6.21%  "
Set-TextValue $ws "D15" "6.01"
Set-TextValue $ws "E15" "  -11.26%  "
Set-TextValue $ws "D16" "13.93"
Set-TextValue $ws "E16" "  -10.73%  "
Set-TextValue $ws "D17" "2.265.72"
Set-TextValue $ws "E17" "  -6.12%  "
Set-TextValue $ws "D18" "0.709"
Set-TextValue $ws "E18" "  -8.07%  "
Set-TextValue $ws "D19" "38.964.90"
Set-TextValue $ws "E19" "  -4.76%  "
Set-TextValue $ws "D20" "0.0₃0852"
Set-TextValue $ws "E20" "  -7.38%  "
Set-TextValue $ws "D21" "5.72"
Set-TextValue $ws "E21" "  -8.15%  "
Set-TextValue $ws "D22" "64.78"
Set-TextValue $ws "E22" "  -8.73%  "
Set-TextValue $ws "D23" "9.82"
Set-TextValue $ws "E23" "  -9.72%  "
Set-TextValue $ws "D24" "225.13"
Set-TextValue $ws "E24" "  -6.34%  "
Set-TextValue $ws "D25" "0.999"
Set-TextValue $ws "E25" "  -0.21%  "
Set-TextValue $ws "D26" "2.37"
Set-TextValue $ws "E26" "  -11.24%  "
Set-TextValue $ws "D27" "1.72"
Set-TextValue $ws "E27" "  -6.75%  "
Set-TextValue $ws "B28" "Toncoin"
Set-TextValue $ws "C28" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws "D28" "2.18"
Set-TextValue $ws "E28" "  -2.42%  "
Set-TextValue $ws "B29" "EthereumClassic"
Set-TextValue $ws "C29" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D29" "22.07"
Set-TextValue $ws "E29" "  -8.24%  "
Set-TextValue $ws "D30" "8.76"
Set-TextValue $ws "E30" "  -8.06%  "
Set-TextValue $ws "D31" "147.66"
Set-TextValue $ws "E31" "  -6.37%  "
Set-TextValue $ws "D32" "31.05"
Set-TextValue $ws "E32" "  -9.38%  "
Set-TextValue $ws "E33" "  -0.01%  "
Set-TextValue $ws "B34" "WEMIXToken"
Set-TextValue $ws "C34" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws "D34" "2.36"
Set-TextValue $ws "E34" "  -4.07%  "
Set-TextValue $ws "B35" "Filecoin"
Set-TextValue $ws "C35" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D35" "4.71"
Set-TextValue $ws "E35" "  -10.70%  "
Set-TextValue $ws "D36" "0.0677"
Set-TextValue $ws "E36" "  -8.26%  "
Set-TextValue $ws "D37" "0.109"
Set-TextValue $ws "E37" "  -4.81%  "
Set-TextValue $ws "B38" "Kaspa"
Set-TextValue $ws "C38" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D38" "0.0943"
Set-TextValue $ws "E38" "  -5.80%  "
Set-TextValue $ws "B39" "LidoDAOToken"
Set-TextValue $ws "C39" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws "D39" "2.60"
Set-TextValue $ws "E39" "  -9.66%  "
Set-TextValue $ws "D40" "14.52"
Set-TextValue $ws "E40" "  -10.83%  "
Set-TextValue $ws "D41" "1.58"
Set-TextValue $ws "E41" "  -10.45%  "
Set-TextValue $ws "D42" "3.58"
Set-TextValue $ws "E42" "  -7.43%  "
Set-TextValue $ws "B43" "ApeXProtocol"
Set-TextValue $ws "C43" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws "D43" "2.22"
Set-TextValue $ws "E43" "  -3.48%  "
Set-TextValue $ws "B44" "Maker"
Set-TextValue $ws "C44" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws "D44" "1.922.41"
Set-TextValue $ws "E44" "  -3.41%  "
Set-TextValue $ws "D45" "0.0253"
Set-TextValue $ws "E45" "  -7.69%  "
Set-TextValue $ws "D46" "16.08"
Set-TextValue $ws "E46" "  -12.45%  "
Set-TextValue $ws "D47" "8.97"
Set-TextValue $ws "E47" "  -4.43%  "
Set-TextValue $ws "D48" "2.50"
Set-TextValue $ws "E48" "  -12.65%  "
Set-TextValue $ws "D49" "2.489.18"
Set-TextValue $ws "E49" "  -5.96%  "
Set-TextValue $ws "D50" "87.13"
Set-TextValue $ws "E50" "  -7.38%  "
Set-TextValue $ws "D51" "65.72"
Set-TextValue $ws "E51" "  -10.98%  "
